$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text corrections (sharedStrings content updates) ---

# Column H (contact info) for every data row: placeholder -> real name
$ws.Range("H4:H14").Value = "Иванов Иван Иванович"

# Fix the "Мустаевский" institution name spacing/hyphenation (rows 6-9, column B)
$ws.Range("B6:B9").Value = "ГБУСО «Мустаевский психоневрологический интернат»"

# Fix the "Сакмарский" institution name spacing/hyphenation for most rows (10,11,13,14)
$ws.Range("B10:B11").Value = "ГБУСО «Сакмарский психоневрологический интернат»"
$ws.Range("B13:B14").Value = "ГБУСО «Сакмарский психоневрологический интернат»"

# Row 12 gets a distinct corrected variant (kept hyphenated "психо-неврологический")
$ws.Range("B12").Value = "ГБУСО «Сакмарский психо-неврологический интернат»"

# --- Row height adjustments (re-layout after the text edits) ---
$ws.Rows.Item(4).RowHeight = 49.25
$ws.Rows.Item(6).RowHeight = 37.3
$ws.Rows.Item(7).RowHeight = 37.3
$ws.Rows.Item(8).RowHeight = 37.3
$ws.Rows.Item(9).RowHeight = 37.3
$ws.Rows.Item(10).RowHeight = 37.3
$ws.Rows.Item(11).RowHeight = 37.3
$ws.Rows.Item(12).RowHeight = 37.3
$ws.Rows.Item(13).RowHeight = 37.3
$ws.Rows.Item(14).RowHeight = 37.3

# --- Selection / view state ---
$ws.Range("I19").Select() | Out-Null
